$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.632.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.489.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "493.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.85%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.494.22"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.76"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0986"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.73%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.916.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.689.96"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.497.76"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.81"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.84"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.410"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.606.76"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.60"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0811"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.37"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.27"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.865"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.39"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0563"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.615"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.95"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "265.32"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.73%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.98%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.23"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.78"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.888.97"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.82%  "
